$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30 data
$dateText = "2025/12/09"
$gameText = "逃离鸭科夫"
$modCount = 1351

# --- A30 (date-looking text must stay as literal text, not be parsed as a date) ---
$helper = $ws.Cells.Item(1,5)
$helper.NumberFormat = "@"
$helper.Value = $dateText

$a30 = $ws.Cells.Item(30,1)
$a30.HorizontalAlignment = -4108
$a30.VerticalAlignment = -4108

$helper.Copy()
$a30.PasteSpecial(-4163)
$helper.Clear()

# --- B30 (plain text) ---
$b30 = $ws.Cells.Item(30,2)
$b30.HorizontalAlignment = -4108
$b30.VerticalAlignment = -4108
$b30.Value = $gameText

# --- C30 (number) ---
$c30 = $ws.Cells.Item(30,3)
$c30.HorizontalAlignment = -4108
$c30.VerticalAlignment = -4108
$c30.Value = $modCount
